# Updated symbol list (coin prices / ranking shuffle) for cryptos.xlsx
# Prices are stored as text in the sheet, so numeric-looking values are
# written with a leading apostrophe to force text entry (avoids Excel
# auto-converting them to numbers and losing formatting like trailing
# zeros / leading zeros / very small magnitudes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'244.52"

# Row 3
$ws.Range("D3").Value = "'23.84"

# Row 5
$ws.Range("D5").Value = "'5.244"

# Row 6
$ws.Range("D6").Value = "'0.05885"

# Row 7
$ws.Range("D7").Value = "'6.475"

# Row 8
$ws.Range("D8").Value = "'3.344"

# Row 9
$ws.Range("D9").Value = "'0.8168"

# Row 10
$ws.Range("D10").Value = "'0.8939"

# Row 11
$ws.Range("B11").Value = 'One'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D11").Value = "'0.0006035"
$ws.Range("E11").Value = '10OneONEWorstin24h'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = "'0.1379"
$ws.Range("E12").Value = '11WazirXWRX'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = "'0.07232"
$ws.Range("E13").Value = '12MandalaExchangeTokenMDX'

# Row 14
$ws.Range("B14").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C14").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D14").Value = "'0.03065"
$ws.Range("E14").Value = '13LiechtensteinCryptoassetsExchangeLCX'

# Row 15
$ws.Range("B15").Value = 'BitrueCoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D15").Value = "'0.03033"
$ws.Range("E15").Value = '14BitrueCoinBTR'

# Row 16
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = "'0.09346"
$ws.Range("E16").Value = '15BitMartTokenBMX'

# Row 17
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = "'3.883"
$ws.Range("E17").Value = '16MCDexMCB'

# Row 18
$ws.Range("B18").Value = 'BitForexToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D18").Value = "'0.001536"
$ws.Range("E18").Value = '17BitForexTokenBF'

# Row 19
$ws.Range("B19").Value = 'CoinExToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D19").Value = "'0.04707"
$ws.Range("E19").Value = '18CoinExTokenCET'

# Row 20
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = "'0.006234"
$ws.Range("E20").Value = '19TigerCashTCH'

# Row 21
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = "'0.001261"
$ws.Range("E21").Value = '20BitKanKAN'

# Row 22
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").Value = "'0.004604"
$ws.Range("E22").Value = '21HotbitTokenHTB'

# Row 23
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = "'0.00008698"
$ws.Range("E23").Value = '22NitroExNTX'

# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = "'2.178"
$ws.Range("E24").Value = '23BTSETokenBTSE'

# Row 26
$ws.Range("D26").Value = "'0.1309"

# Row 28
$ws.Range("D28").Value = "'0.0002338"

# Row 40
$ws.Range("D40").Value = "'0.03791"

# Row 41
$ws.Range("D41").Value = "'0.006368"

# Row 42
$ws.Range("D42").Value = "'0.1058"

# Row 43
$ws.Range("D43").Value = "'0.002606"

# Row 44
$ws.Range("D44").Value = "'0.007071"

# Row 45
$ws.Range("D45").Value = "'0.00005392"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"

# Row 47
$ws.Range("D47").Value = "'0.5397"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

# Row 48
$ws.Range("D48").Value = "'0.02105"
$ws.Range("E48").Value = '47BOLOBOLOBestin24h'

# Row 49
$ws.Range("D49").Value = "'0.00002099"

# Row 50
$ws.Range("D50").Value = "'0.0001999"
